# Fixed issues with 81RF protective element
# Changed default xls parameters to disable 81x protections.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("relays")

# Disable 81RF protection by updating its default pickup/delay settings
# for the three relay rows (2-4), in both the primary (T:V) and
# secondary/backup (AH:AJ) 81RF column groups.
foreach ($r in 2..4) {
    $ws.Range("T$r").Value = 100
    $ws.Range("U$r").Value = 10
    $ws.Range("V$r").Value = 0.1

    $ws.Range("AH$r").Value = 100
    $ws.Range("AI$r").Value = 10
    $ws.Range("AJ$r").Value = 0.1
}

# Reselect the updated 81RF cell range to mirror the edited selection.
$ws.Range("AH2:AJ4").Select()
